# Turn the run "Hi " into three separate runs "H" | "ey" | " " (i.e. the
# word "Hi" becomes "Hey"), while leaving the following "World" / "!"
# runs byte-for-byte untouched (same text, same rPr, same w:rsidR).
#
# A direct text-replace (Range.Text / Range.InsertAfter / Range.Delete)
# coalesces *every* same-formatted run in the paragraph into a single run,
# which would also swallow "World" and "!" into the edited run and drop
# their w:rsidR. To avoid that we temporarily "shield" the "World" run
# with a throw-away direct-formatting toggle (Bold on/off) so the
# text-replace's run-coalescing stops at the shield boundary, then we use
# the same on/off formatting-toggle trick (which only splits runs at the
# touched boundaries, it does not coalesce the whole paragraph) to carve
# the edited prefix into the three target runs.

$d = $word.ActiveDocument

# Locate the literal "Hi " run and the following "World" word dynamically
# so this does not depend on hard-coded character offsets.
$hi = $d.Content
$hi.Find.Execute("Hi ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null

$hiStart = $hi.Start          # start of "Hi "
$hiEnd = $hi.End              # end of "Hi " (== start of "World")

$world = $d.Content
$world.Find.Execute("World", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null

$worldStart = $world.Start
$worldEnd = $world.End

# 1) Shield "World" with a no-op direct-formatting toggle so the upcoming
#    text edit's run-coalescing cannot merge into (or swallow the
#    w:rsidR off) the "World" run.
$shield = $d.Range($worldStart, $worldEnd)
$shield.Bold = $true

# 2) Replace "i" with "ey" (the only actual text content change).
$iStart = $hiStart + 1
$iEnd = $iStart + 1
$iRange = $d.Range($iStart, $iEnd)
$iRange.Text = "ey"

# Text grew by 1 character ("i" -> "ey"), so everything from here on is
# shifted by +1 relative to the original offsets.
$shiftedWorldStart = $worldStart + 1
$shiftedWorldEnd = $worldEnd + 1

# 3) Remove the shield.
$unshield = $d.Range($shiftedWorldStart, $shiftedWorldEnd)
$unshield.Bold = $false

# 4) Split the (now merged) "H", "ey", " " prefix apart using the same
#    formatting-toggle trick -- each toggle only splits runs at the
#    touched boundary, it does not re-coalesce the paragraph.
$spaceStart = $shiftedWorldStart - 1
$spaceEnd = $shiftedWorldStart

$spaceRange = $d.Range($spaceStart, $spaceEnd)
$spaceRange.Bold = $true
$spaceRange.Bold = $false

$eyRange = $d.Range($iStart, $spaceStart)
$eyRange.Bold = $true
$eyRange.Bold = $false
